# Apply "Copied and renamed classes" edit to StateTransitionTable.xlsx
#
# The state-transition labels "Ideal Temp" and "Timer @ 10" that were used
# as header text in the second table (row 8, columns F and I) were copied
# into their own distinct strings and renamed:
#   F8: "Ideal Temp"  -> "Ideal Temperature"
#   I8: "Timer @ 10"  -> "TimerAt10"
# Every other cell that used to share a string index with these two values
# keeps its original text (it still references the original shared string,
# which is otherwise untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F8").Value = "Ideal Temperature"
$ws.Range("I8").Value = "TimerAt10"

# Restore the view state captured in the saved workbook: scrolled so column
# E is the left-most visible column, with I8 selected as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I8").Select()
